# Applies the Behemoth_Profits price-refresh update described by the commit diff.
# 216 numeric cell updates + 1 cell clear (M129 on CUL, which the diff removes
# outright rather than leaving a stale value) across all 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 554.73334
$ws.Cells.Item(28, 9).Value = 378.9524
$ws.Cells.Item(28, 10).Value = 964.8889
$ws.Cells.Item(28, 11).Value = 378.9524
$ws.Cells.Item(28, 12).Value = 964.8889
$ws.Cells.Item(28, 13).Value = 106.0476
$ws.Cells.Item(28, 14).Value = -1934.8889
$ws.Cells.Item(64, 8).Value = 4631.5264
$ws.Cells.Item(64, 9).Value = 3599.8
$ws.Cells.Item(64, 10).Value = 5000
$ws.Cells.Item(64, 11).Value = 3599.8
$ws.Cells.Item(64, 12).Value = 5000
$ws.Cells.Item(64, 13).Value = -3351.8
$ws.Cells.Item(64, 14).Value = -5496
$ws.Cells.Item(67, 8).Value = 4631.5264
$ws.Cells.Item(67, 9).Value = 3599.8
$ws.Cells.Item(67, 10).Value = 5000
$ws.Cells.Item(67, 11).Value = 3599.8
$ws.Cells.Item(67, 12).Value = 5000
$ws.Cells.Item(67, 13).Value = -2741.8
$ws.Cells.Item(67, 14).Value = -6716
$ws.Cells.Item(74, 8).Value = 3680.2
$ws.Cells.Item(74, 9).Value = 3725.25
$ws.Cells.Item(74, 11).Value = 3725.25
$ws.Cells.Item(74, 13).Value = -2789.25
$ws.Cells.Item(76, 8).Value = 6140.3706
$ws.Cells.Item(76, 9).Value = 5988.8887
$ws.Cells.Item(76, 11).Value = 5988.8887
$ws.Cells.Item(76, 13).Value = -5673.8887
$ws.Cells.Item(77, 8).Value = 3680.2
$ws.Cells.Item(77, 9).Value = 3725.25
$ws.Cells.Item(77, 11).Value = 18626.25
$ws.Cells.Item(77, 13).Value = -13946.25
$ws.Cells.Item(79, 8).Value = 6140.3706
$ws.Cells.Item(79, 9).Value = 5988.8887
$ws.Cells.Item(79, 11).Value = 5988.8887
$ws.Cells.Item(79, 13).Value = -4896.8887
$ws.Cells.Item(112, 8).Value = 2074.0667
$ws.Cells.Item(112, 10).Value = 2143.7144
$ws.Cells.Item(112, 12).Value = 6431.1432
$ws.Cells.Item(112, 14).Value = -8647.143199999999
$ws.Cells.Item(132, 8).Value = 1327.2
$ws.Cells.Item(132, 9).Value = 773.0294
$ws.Cells.Item(132, 11).Value = 2319.0882
$ws.Cells.Item(132, 13).Value = 210.9117999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 25002434
$ws.Cells.Item(45, 9).Value = 38463292
$ws.Cells.Item(45, 11).Value = 38463292
$ws.Cells.Item(45, 13).Value = -38462915
$ws.Cells.Item(60, 8).Value = 40000
$ws.Cells.Item(60, 10).Value = 40000
$ws.Cells.Item(60, 12).Value = 40000
$ws.Cells.Item(60, 14).Value = -41466
$ws.Cells.Item(132, 8).Value = 3537.8845
$ws.Cells.Item(132, 9).Value = 3537.8845
$ws.Cells.Item(132, 11).Value = 10613.6535
$ws.Cells.Item(132, 13).Value = -8083.6535
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 599.5
$ws.Cells.Item(5, 10).Value = 105
$ws.Cells.Item(5, 12).Value = 105
$ws.Cells.Item(5, 14).Value = -331
$ws.Cells.Item(86, 8).Value = 2336.5386
$ws.Cells.Item(86, 9).Value = 2688.1667
$ws.Cells.Item(86, 11).Value = 2688.1667
$ws.Cells.Item(86, 13).Value = -1565.1667
$ws.Cells.Item(89, 8).Value = 2336.5386
$ws.Cells.Item(89, 9).Value = 2688.1667
$ws.Cells.Item(89, 11).Value = 13440.8335
$ws.Cells.Item(89, 13).Value = -7824.833500000001
$ws.Cells.Item(92, 8).Value = 66326.664
$ws.Cells.Item(92, 10).Value = 66326.664
$ws.Cells.Item(92, 12).Value = 66326.664
$ws.Cells.Item(92, 14).Value = -71318.664
$ws.Cells.Item(105, 8).Value = 2085.125
$ws.Cells.Item(105, 9).Value = 1034.6666
$ws.Cells.Item(105, 10).Value = 2715.4
$ws.Cells.Item(105, 11).Value = 1034.6666
$ws.Cells.Item(105, 12).Value = 2715.4
$ws.Cells.Item(105, 13).Value = 712.3334
$ws.Cells.Item(105, 14).Value = -6209.4
$ws.Cells.Item(134, 8).Value = 787785.9
$ws.Cells.Item(134, 9).Value = 1832.3334
$ws.Cells.Item(134, 11).Value = 5497.0002
$ws.Cells.Item(134, 13).Value = -2962.0002
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 5600
$ws.Cells.Item(25, 9).Value = 2400
$ws.Cells.Item(25, 10).Value = 12000
$ws.Cells.Item(25, 11).Value = 2400
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = -2226
$ws.Cells.Item(25, 14).Value = -12348
$ws.Cells.Item(31, 8).Value = 806415.1
$ws.Cells.Item(31, 9).Value = 2377.85
$ws.Cells.Item(31, 10).Value = 2593164.5
$ws.Cells.Item(31, 11).Value = 2377.85
$ws.Cells.Item(31, 12).Value = 2593164.5
$ws.Cells.Item(31, 13).Value = -2082.85
$ws.Cells.Item(31, 14).Value = -2593754.5
$ws.Cells.Item(34, 8).Value = 806415.1
$ws.Cells.Item(34, 9).Value = 2377.85
$ws.Cells.Item(34, 10).Value = 2593164.5
$ws.Cells.Item(34, 11).Value = 2377.85
$ws.Cells.Item(34, 12).Value = 2593164.5
$ws.Cells.Item(34, 13).Value = -2175.85
$ws.Cells.Item(34, 14).Value = -2593568.5
$ws.Cells.Item(62, 8).Value = 1672982.4
$ws.Cells.Item(62, 9).Value = 1672982.4
$ws.Cells.Item(62, 11).Value = 1672982.4
$ws.Cells.Item(62, 13).Value = -1672358.4
$ws.Cells.Item(65, 8).Value = 1672982.4
$ws.Cells.Item(65, 9).Value = 1672982.4
$ws.Cells.Item(65, 11).Value = 8364912
$ws.Cells.Item(65, 13).Value = -8361792
$ws.Cells.Item(134, 8).Value = 3062.1853
$ws.Cells.Item(134, 9).Value = 2377.4783
$ws.Cells.Item(134, 11).Value = 7132.4349
$ws.Cells.Item(134, 13).Value = -4597.4349
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 114.46667
$ws.Cells.Item(38, 9).Value = 120.9
$ws.Cells.Item(38, 11).Value = 362.7
$ws.Cells.Item(38, 13).Value = -15.70000000000005
$ws.Cells.Item(129, 8).Value = 41671920
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 41671920
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 125015760
$ws.Cells.Item(129, 14).Value = -125025760
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(131, 8).Value = 153530.89
$ws.Cells.Item(131, 10).Value = 11537.0625
$ws.Cells.Item(131, 12).Value = 34611.1875
$ws.Cells.Item(131, 14).Value = -44691.1875
$ws.Cells.Item(132, 8).Value = 1934.4445
$ws.Cells.Item(132, 9).Value = 1303.3334
$ws.Cells.Item(132, 11).Value = 11730.0006
$ws.Cells.Item(132, 13).Value = -9200.000599999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9205.772000000001
$ws.Cells.Item(70, 9).Value = 10158
$ws.Cells.Item(70, 10).Value = 6666.5
$ws.Cells.Item(70, 11).Value = 10158
$ws.Cells.Item(70, 12).Value = 6666.5
$ws.Cells.Item(70, 13).Value = -9888
$ws.Cells.Item(70, 14).Value = -7206.5
$ws.Cells.Item(73, 8).Value = 9205.772000000001
$ws.Cells.Item(73, 9).Value = 10158
$ws.Cells.Item(73, 10).Value = 6666.5
$ws.Cells.Item(73, 11).Value = 10158
$ws.Cells.Item(73, 12).Value = 6666.5
$ws.Cells.Item(73, 13).Value = -9222
$ws.Cells.Item(73, 14).Value = -8538.5
$ws.Cells.Item(75, 8).Value = 29899.1
$ws.Cells.Item(75, 10).Value = 29899.1
$ws.Cells.Item(75, 12).Value = 29899.1
$ws.Cells.Item(75, 14).Value = -31647.1
$ws.Cells.Item(78, 8).Value = 29899.1
$ws.Cells.Item(78, 10).Value = 29899.1
$ws.Cells.Item(78, 12).Value = 89697.29999999999
$ws.Cells.Item(78, 14).Value = -98433.29999999999
$ws.Cells.Item(102, 8).Value = 2255.9412
$ws.Cells.Item(102, 9).Value = 1668.64
$ws.Cells.Item(102, 11).Value = 1668.64
$ws.Cells.Item(102, 13).Value = -46.6400000000001
$ws.Cells.Item(132, 8).Value = 200020530
$ws.Cells.Item(132, 9).Value = 250000660
$ws.Cells.Item(132, 11).Value = 750001980
$ws.Cells.Item(132, 13).Value = -749999450
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 12499.833
$ws.Cells.Item(13, 9).Value = 13749.75
$ws.Cells.Item(13, 11).Value = 13749.75
$ws.Cells.Item(13, 13).Value = -13609.75
$ws.Cells.Item(46, 8).Value = 2598.8572
$ws.Cells.Item(46, 9).Value = 2178.4
$ws.Cells.Item(46, 10).Value = 3650
$ws.Cells.Item(46, 11).Value = 2178.4
$ws.Cells.Item(46, 12).Value = 3650
$ws.Cells.Item(46, 13).Value = -1990.4
$ws.Cells.Item(46, 14).Value = -4026
$ws.Cells.Item(82, 8).Value = 1657
$ws.Cells.Item(82, 9).Value = 1775
$ws.Cells.Item(82, 10).Value = 1609.8
$ws.Cells.Item(82, 11).Value = 1775
$ws.Cells.Item(82, 12).Value = 1609.8
$ws.Cells.Item(82, 13).Value = -1414
$ws.Cells.Item(82, 14).Value = -2331.8
$ws.Cells.Item(85, 8).Value = 1657
$ws.Cells.Item(85, 9).Value = 1775
$ws.Cells.Item(85, 10).Value = 1609.8
$ws.Cells.Item(85, 11).Value = 1775
$ws.Cells.Item(85, 12).Value = 1609.8
$ws.Cells.Item(85, 13).Value = -527
$ws.Cells.Item(85, 14).Value = -4105.8
$ws.Cells.Item(122, 8).Value = 5858.6
$ws.Cells.Item(122, 9).Value = 5551.0527
$ws.Cells.Item(122, 11).Value = 16653.1581
$ws.Cells.Item(122, 13).Value = -14203.1581
$ws.Cells.Item(132, 8).Value = 918934.5600000001
$ws.Cells.Item(132, 9).Value = 102521
$ws.Cells.Item(132, 10).Value = 5001002.5
$ws.Cells.Item(132, 11).Value = 307563
$ws.Cells.Item(132, 12).Value = 15003007.5
$ws.Cells.Item(132, 13).Value = -305033
$ws.Cells.Item(132, 14).Value = -15008067.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2898.3438
$ws.Cells.Item(122, 9).Value = 2916.2273
$ws.Cells.Item(122, 11).Value = 8748.6819
$ws.Cells.Item(122, 13).Value = -6298.6819
$ws.Cells.Item(132, 8).Value = 403719.44
$ws.Cells.Item(132, 9).Value = 3866.2354
$ws.Cells.Item(132, 11).Value = 11598.7062
$ws.Cells.Item(132, 13).Value = -9068.706200000001
$ws.Cells.Item(136, 8).Value = 6393.9375
$ws.Cells.Item(136, 9).Value = 6432.8
$ws.Cells.Item(136, 10).Value = 6329.1665
$ws.Cells.Item(136, 11).Value = 19298.4
$ws.Cells.Item(136, 12).Value = 18987.4995
$ws.Cells.Item(136, 13).Value = -16748.4
$ws.Cells.Item(136, 14).Value = -24087.4995

Write-Host "Applied 216 cell updates and 1 cell clear across 8 sheets."